# Lab5A instructions: small text edits.
$d = $word.ActiveDocument

# 1) "Part 2: Monitoring HTTP Traffic and Debugging" - collapse the four
#    runs ("Part ", "2", ": ", "Monitoring HTTP Traffic and Debugging")
#    into one by doing a literal find/replace over the whole phrase.
$d.Content.Find.Execute("Part 2: Monitoring HTTP Traffic and Debugging", `
    $false, $true, $false, $false, $false, $true, 1, $false, `
    "Part 2: Monitoring HTTP Traffic and Debugging", 2)

# 2) "Submit the following:" - collapse the two runs
#    ("Submit the following", ":") into one.
$d.Content.Find.Execute("Submit the following:", `
    $false, $true, $false, $false, $false, $true, 1, $false, `
    "Submit the following:", 2)

# 3) Header (section 1's "Group A - Deploy to Azure" title): append
#    " + Debugging" after "Deploy to Azure".
foreach ($sec in $d.Sections) {
    $hdr = $sec.Headers.Item(2)
    if ($hdr.Exists) {
        $hdr.Range.Find.Execute("Group A " + [char]0x2013 + " Deploy to Azure", `
            $false, $true, $false, $false, $false, $true, 1, $false, `
            "Group A " + [char]0x2013 + " Deploy to Azure + Debugging", 2)
    }
}
